$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 4589.1816
$ws.Range("J19").Value = 6784.143
$ws.Range("L19").Value = 6784.143
$ws.Range("N19").Value = -7134.143
# row 88
$ws.Range("H88").Value = 2552.5151
$ws.Range("I88").Value = 1788.75
$ws.Range("J88").Value = 2657.862
$ws.Range("K88").Value = 1788.75
$ws.Range("L88").Value = 2657.862
$ws.Range("M88").Value = -1382.75
$ws.Range("N88").Value = -3469.862
# row 91
$ws.Range("H91").Value = 2552.5151
$ws.Range("I91").Value = 1788.75
$ws.Range("J91").Value = 2657.862
$ws.Range("K91").Value = 1788.75
$ws.Range("L91").Value = 2657.862
$ws.Range("M91").Value = -384.75
$ws.Range("N91").Value = -5465.862
# row 98
$ws.Range("H98").Value = 1714.1364
$ws.Range("I98").Value = 1652.1666
$ws.Range("J98").Value = 1993
$ws.Range("K98").Value = 1652.1666
$ws.Range("L98").Value = 1993
$ws.Range("M98").Value = -154.1666
$ws.Range("N98").Value = -4989
# row 100
$ws.Range("H100").Value = 2478.8572
$ws.Range("I100").Value = 2808.6667
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 2808.6667
$ws.Range("L100").Value = 500
$ws.Range("M100").Value = -2267.6667
$ws.Range("N100").Value = -1582
# row 122
$ws.Range("H122").Value = 1714.1364
$ws.Range("I122").Value = 1652.1666
$ws.Range("J122").Value = 1993
$ws.Range("K122").Value = 4956.4998
$ws.Range("L122").Value = 5979
$ws.Range("M122").Value = -2506.4998
$ws.Range("N122").Value = -10879
# row 132
$ws.Range("H132").Value = 28575530
$ws.Range("I132").Value = 30307320
$ws.Range("J132").Value = 975
$ws.Range("K132").Value = 90921960
$ws.Range("L132").Value = 2925
$ws.Range("M132").Value = -90919430
$ws.Range("N132").Value = -7985
# row 137
$ws.Range("H137").Value = 54532.47
$ws.Range("I137").Value = 149785.92
$ws.Range("J137").Value = 2576.0454
$ws.Range("K137").Value = 449357.76
$ws.Range("L137").Value = 7728.1362
$ws.Range("M137").Value = -446807.76
$ws.Range("N137").Value = -12828.1362
# row 141
$ws.Range("H141").Value = 2127.2666
$ws.Range("I141").Value = 1863.25
$ws.Range("K141").Value = 5589.75
$ws.Range("M141").Value = -409.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 5154.931
$ws.Range("I32").Value = 3673.8057
$ws.Range("J32").Value = 7578.591
$ws.Range("K32").Value = 3673.8057
$ws.Range("L32").Value = 7578.591
$ws.Range("M32").Value = -3386.8057
$ws.Range("N32").Value = -8152.591
# row 61
$ws.Range("H61").Value = 2421.15
$ws.Range("I61").Value = 1827.4286
$ws.Range("K61").Value = 1827.4286
$ws.Range("M61").Value = -1615.4286
# row 63
$ws.Range("H63").Value = 3761.3635
$ws.Range("I63").Value = 3052.7778
$ws.Range("K63").Value = 3052.7778
$ws.Range("M63").Value = -2366.7778
# row 66
$ws.Range("H66").Value = 3761.3635
$ws.Range("I66").Value = 3052.7778
$ws.Range("K66").Value = 15263.889
$ws.Range("M66").Value = -11831.889
# row 132
$ws.Range("H132").Value = 2980.0625
$ws.Range("I132").Value = 1929
$ws.Range("J132").Value = 3797.5557
$ws.Range("K132").Value = 5787
$ws.Range("L132").Value = 11392.6671
$ws.Range("M132").Value = -3257
$ws.Range("N132").Value = -16452.6671
# row 136
$ws.Range("H136").Value = 2421.15
$ws.Range("I136").Value = 1827.4286
$ws.Range("K136").Value = 5482.2858
$ws.Range("M136").Value = -2932.2858

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 105
$ws.Range("H105").Value = 4167928.5
$ws.Range("I105").Value = 4465352
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 4465352
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -4463605
$ws.Range("N105").Value = -7494

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 441.8
$ws.Range("I7").Value = 143.75
$ws.Range("K7").Value = 143.75
$ws.Range("M7").Value = -30.75
# row 58
$ws.Range("H58").Value = 2966.6572
$ws.Range("I58").Value = 3242.56
$ws.Range("K58").Value = 3242.56
$ws.Range("M58").Value = -3039.56
# row 134
$ws.Range("H134").Value = 3292.0557
$ws.Range("I134").Value = 2632.7144
$ws.Range("K134").Value = 7898.1432
$ws.Range("M134").Value = -5363.1432
# row 136
$ws.Range("H136").Value = 2966.6572
$ws.Range("I136").Value = 3242.56
$ws.Range("K136").Value = 9727.68
$ws.Range("M136").Value = -7177.68
# row 141
$ws.Range("H141").Value = 29860.8
$ws.Range("J141").Value = 29826
$ws.Range("L141").Value = 29826
$ws.Range("N141").Value = -40186

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 68
$ws.Range("H68").Value = 936.25
$ws.Range("I68").Value = 718
$ws.Range("K68").Value = 2154
$ws.Range("M68").Value = -1343
# row 71
$ws.Range("H71").Value = 936.25
$ws.Range("I71").Value = 718
$ws.Range("K71").Value = 6462
$ws.Range("M71").Value = -2406
# row 122
$ws.Range("H122").Value = 1048.0714
$ws.Range("J122").Value = 947.2222
$ws.Range("L122").Value = 8524.9998
$ws.Range("N122").Value = -13424.9998
# row 139
$ws.Range("H139").Value = 2098.5
$ws.Range("I139").Value = 2098.5
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 6295.5
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -1155.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 1880107.4
$ws.Range("I80").Value = 3486621.2
$ws.Range("K80").Value = 3486621.2
$ws.Range("M80").Value = -3485623.2
# row 83
$ws.Range("H83").Value = 1880107.4
$ws.Range("I83").Value = 3486621.2
$ws.Range("K83").Value = 17433106
$ws.Range("M83").Value = -17428114
# row 102
$ws.Range("H102").Value = 5439120
$ws.Range("I102").Value = 6945940
$ws.Range("J102").Value = 2760328.2
$ws.Range("K102").Value = 6945940
$ws.Range("L102").Value = 2760328.2
$ws.Range("M102").Value = -6944318
$ws.Range("N102").Value = -2763572.2
# row 122
$ws.Range("H122").Value = 320002.75
$ws.Range("I122").Value = 469525.3
$ws.Range("J122").Value = 4344
$ws.Range("K122").Value = 1408575.9
$ws.Range("L122").Value = 13032
$ws.Range("M122").Value = -1406125.9
$ws.Range("N122").Value = -17932
# row 132
$ws.Range("H132").Value = 2896.9412
$ws.Range("I132").Value = 2546.7812
$ws.Range("J132").Value = 8499.5
$ws.Range("K132").Value = 7640.3436
$ws.Range("L132").Value = 25498.5
$ws.Range("M132").Value = -5110.3436
$ws.Range("N132").Value = -30558.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 10
$ws.Range("H10").Value = 67958.75
$ws.Range("J10").Value = 199996.25
$ws.Range("L10").Value = 199996.25
$ws.Range("N10").Value = -200276.25
# row 46
$ws.Range("H46").Value = 5934.6523
$ws.Range("I46").Value = 5490
$ws.Range("J46").Value = 6276.6924
$ws.Range("K46").Value = 5490
$ws.Range("L46").Value = 6276.6924
$ws.Range("M46").Value = -5302
$ws.Range("N46").Value = -6652.6924
# row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0
# row 136
$ws.Range("H136").Value = 56202.71
$ws.Range("I136").Value = 77420.63
$ws.Range("K136").Value = 232261.89
$ws.Range("M136").Value = -229711.89
